# Actualizacion automatica hashcode - update hash values in column B
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(11, 2).Value = "1f682c4baf00039722b9d3b2a8f6431f"  # 05-050301A
$ws.Cells.Item(44, 2).Value = "a2cfcbfef9b7b4aed5ed06cdf76e820f"  # 05-050105A
$ws.Cells.Item(74, 2).Value = "9555bf74da8a390313ded720eb47dce7"  # 05-050103A
$ws.Cells.Item(89, 2).Value = "160ee88f449d69ffbf488ebe9d2dcc44"  # 05-050104A
$ws.Cells.Item(99, 2).Value = "ec5bd2a050b8a245967e920be6cdaaa2"  # 05-050101A
$ws.Cells.Item(100, 2).Value = "85819c9b0ee864700a6fb3abf7b62758"  # 04-040021TM
$ws.Cells.Item(104, 2).Value = "afc45b0ea45fcd2114d8102997488408"  # 04-040021TP
$ws.Cells.Item(110, 2).Value = "4050bd447a74401c61ea746f9711d4fc"  # 05-050102A
$ws.Cells.Item(113, 2).Value = "956b266fd844e9f3fca2194ee278fadb"  # 04-040021TC
$ws.Cells.Item(121, 2).Value = "81667d4f5140992663fc6287a415e11f"  # 05-050301TP
$ws.Cells.Item(122, 2).Value = "d15ca3c8fb72fbbd22db7c2394f28a69"  # 04-040014TC
$ws.Cells.Item(123, 2).Value = "3c295a675ead62d8afffd86dda0453ab"  # 05-050301TC
$ws.Cells.Item(164, 2).Value = "0a80cf60deec27272e68c8141fbee685"  # 04-040021A
$ws.Cells.Item(168, 2).Value = "36c8cd53ba8a46717318adc0a51706b1"  # 05-050105TP
$ws.Cells.Item(191, 2).Value = "aec159b771e496e8cb54e48f8a239e8e"  # 05-050314TP
$ws.Cells.Item(230, 2).Value = "a7ccd9496d18261177551264266f67e7"  # 04-040014TP
$ws.Cells.Item(233, 2).Value = "380c5e4c6ed05e85df43317f9a0cfa66"  # 04-040014TM
$ws.Cells.Item(278, 2).Value = "4f4e6e1d7f91885a3a4f184b8ac396e3"  # 05-050101TP
$ws.Cells.Item(331, 2).Value = "d9986ed4380897b50d61c0803314de7c"  # 04-040018TP
$ws.Cells.Item(342, 2).Value = "052d5b4453144717d9154004c40aed09"  # 04-040018TC
$ws.Cells.Item(343, 2).Value = "9c8e173b79f48d63f00af95644862e76"  # 04-040018TM
$ws.Cells.Item(345, 2).Value = "183913fecc02620ae6913e0667b17656"  # 05-050103TP
$ws.Cells.Item(419, 2).Value = "930e9bd628ccd09c643cd2b4a4b8cfad"  # 05-0709-070905BTC
$ws.Cells.Item(480, 2).Value = "1fd9ef0f8869fc52d6c81138b24ec41c"  # 05-050314A
$ws.Cells.Item(515, 2).Value = "ef4292e83e9c1fb6f80576ab1bfe45d3"  # 05-050312A
$ws.Cells.Item(547, 2).Value = "12134a6651c6de21c72dc6c1e1dae89a"  # 05-050201A
$ws.Cells.Item(619, 2).Value = "bd09cfb4e9f5a5a1edc58ee2f6cbef23"  # 04-040015TC
$ws.Cells.Item(623, 2).Value = "5df9e1ffb7ca51b90d6720532ccfee6f"  # 04-040015TP
$ws.Cells.Item(628, 2).Value = "ae8a27b09551a4de674da30e82a0e23c"  # 04-040015TM
$ws.Cells.Item(726, 2).Value = "0d69356325dde6912774de5cc26ca562"  # 05-050315TC
$ws.Cells.Item(733, 2).Value = "41a70b09bf76f235b51a465777177226"  # 05-050207A
$ws.Cells.Item(768, 2).Value = "8a866f38cea4d509d812189b47eef642"  # 05-050102TP
$ws.Cells.Item(779, 2).Value = "babf3fd530aff2ea45435a4292853ff1"  # 04-040018A
$ws.Cells.Item(816, 2).Value = "1951623ae9020a139ec3467817acc2ab"  # 05-050104TC
$ws.Cells.Item(818, 2).Value = "4c2ed9e49577e877cba8646fab52dc00"  # 04-040015A
$ws.Cells.Item(827, 2).Value = "fe391b223dd9b3e7fc6a5f6ebd9890a3"  # 05-050104TP
$ws.Cells.Item(831, 2).Value = "3ebef27ff7385eb5bb0c6c1d9dc07834"  # 04-040014A
$ws.Cells.Item(874, 2).Value = "c9c849f03081bb7a17b5eba5feebb7ea"  # 03-030032A
$ws.Cells.Item(904, 2).Value = "162cf7f74e2a908d24fead2084dcf5fb"  # 05-050001A
